# Apply the updated cryptocurrency price/volume snapshot to Sheet1.
# D-column values that look like plain numbers are written with a leading
# apostrophe (text-entry marker) so Excel stores them as text, matching
# the original inlineStr cells (e.g. "7.70" must stay "7.70", not become 7.7).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.528.60'
$ws.Range("E2").Value = '  -5.52%  '
# Row 3
$ws.Range("D3").Value = '2.591.28'
$ws.Range("E3").Value = '  -0.16%  '
# Row 4
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.08%  '
# Row 5
$ws.Range("D5").Value = '''299.84'
$ws.Range("E5").Value = '  -2.42%  '
# Row 6
$ws.Range("D6").Value = '''95.33'
$ws.Range("E6").Value = '  -4.15%  '
# Row 7
$ws.Range("D7").Value = '''0.572'
$ws.Range("E7").Value = '  -4.74%  '
# Row 8
$ws.Range("E8").Value = '  -0.07%  '
# Row 9
$ws.Range("D9").Value = '''0.547'
$ws.Range("E9").Value = '  -5.41%  '
# Row 10
$ws.Range("D10").Value = '''36.58'
$ws.Range("E10").Value = '  -6.53%  '
# Row 11
$ws.Range("D11").Value = '''0.0805'
$ws.Range("E11").Value = '  -4.58%  '
# Row 12
$ws.Range("D12").Value = '''7.70'
$ws.Range("E12").Value = '  -5.84%  '
# Row 13
$ws.Range("D13").Value = '2.996.74'
$ws.Range("E13").Value = '  +0.03%  '
# Row 14
$ws.Range("D14").Value = '''0.106'
$ws.Range("E14").Value = '  +0.73%  '
# Row 15
$ws.Range("D15").Value = '2.617.11'
$ws.Range("E15").Value = '  +0.08%  '
# Row 16
$ws.Range("D16").Value = '''0.877'
$ws.Range("E16").Value = '  -5.10%  '
# Row 17
$ws.Range("D17").Value = '''14.16'
$ws.Range("E17").Value = '  -5.53%  '
# Row 18
$ws.Range("D18").Value = '43.478.02'
$ws.Range("E18").Value = '  -6.07%  '
# Row 19
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = "0.0$([char]0x2083)0962"
$ws.Range("E19").Value = '  -4.75%  '
# Row 20
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '''6.54'
$ws.Range("E20").Value = '  -2.74%  '
# Row 21
$ws.Range("D21").Value = '''12.17'
$ws.Range("E21").Value = '  -6.38%  '
# Row 22
$ws.Range("D22").Value = '''72.69'
$ws.Range("E22").Value = '  +1.92%  '
# Row 23
$ws.Range("D23").Value = '''262.48'
$ws.Range("E23").Value = '  -3.69%  '
# Row 24
$ws.Range("D24").Value = '''2.91'
$ws.Range("E24").Value = '  -4.00%  '
# Row 25
$ws.Range("D25").Value = '''2.18'
$ws.Range("E25").Value = '  +0.48%  '
# Row 26
$ws.Range("D26").Value = '''28.91'
$ws.Range("E26").Value = '  -3.12%  '
# Row 27
$ws.Range("E27").Value = '  +0.13%  '
# Row 28
$ws.Range("D28").Value = '''10.08'
$ws.Range("E28").Value = '  -4.68%  '
# Row 29
$ws.Range("D29").Value = '''2.21'
$ws.Range("E29").Value = '  -4.43%  '
# Row 30
$ws.Range("D30").Value = '''37.15'
$ws.Range("E30").Value = '  -4.44%  '
# Row 31
$ws.Range("D31").Value = '''5.88'
$ws.Range("E31").Value = '  -5.92%  '
# Row 32
$ws.Range("D32").Value = '''3.55'
$ws.Range("E32").Value = '  -2.25%  '
# Row 33
$ws.Range("D33").Value = '''2.19'
$ws.Range("E33").Value = '  +0.07%  '
# Row 34
$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").Value = '''2.79'
$ws.Range("E34").Value = '  -1.63%  '
# Row 35
$ws.Range("B35").Value = 'Monero'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D35").Value = '''151.06'
$ws.Range("E35").Value = '  +0.71%  '
# Row 36
$ws.Range("D36").Value = '''0.0798'
$ws.Range("E36").Value = '  -4.89%  '
# Row 37
$ws.Range("D37").Value = '''0.116'
$ws.Range("E37").Value = '  -4.33%  '
# Row 38
$ws.Range("D38").Value = '''24.32'
$ws.Range("E38").Value = '  +5.08%  '
# Row 39
$ws.Range("D39").Value = '''0.118'
$ws.Range("E39").Value = '  -2.84%  '
# Row 40
$ws.Range("D40").Value = '''16.29'
$ws.Range("E40").Value = '  +2.78%  '
# Row 41
$ws.Range("D41").Value = '''3.41'
$ws.Range("E41").Value = '  -5.72%  '
# Row 42
$ws.Range("D42").Value = '''0.0310'
$ws.Range("E42").Value = '  -5.58%  '
# Row 43
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '2.069.17'
$ws.Range("E43").Value = '  -4.53%  '
# Row 44
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '''3.78'
$ws.Range("E44").Value = '  -7.40%  '
# Row 45
$ws.Range("E45").Value = '  -0.09%  '
# Row 46
$ws.Range("D46").Value = '''87.42'
$ws.Range("E46").Value = '  -6.70%  '
# Row 47
$ws.Range("D47").Value = '''9.11'
$ws.Range("E47").Value = '  -4.82%  '
# Row 48
$ws.Range("D48").Value = '2.856.13'
$ws.Range("E48").Value = '  +0.31%  '
# Row 49
$ws.Range("E49").Value = '  +3.22%  '
# Row 50
$ws.Range("D50").Value = '''104.84'
$ws.Range("E50").Value = '  -3.89%  '
# Row 51
$ws.Range("D51").Value = '''0.188'
$ws.Range("E51").Value = '  -5.63%  '
